# Applies the "Add new NCC with feature" edit:
#   - ACR sheet: replace the uniform 0.5 grid (rows 1-5, cols A:CW) with a
#     two-level step pattern (10/43 then 33/43), the switch-over column
#     varying per row.
#   - FAR sheet: three single-cell corrections (AY1, AY3: 0->1; AX4: 1->0)
#   - FRR sheet: the complementary three single-cell corrections
#     (AY1, AY3: 1->0; AX4: 0->1)

$wb = $excel.ActiveWorkbook

$lowValue  = 10 / 43    # 0.23255813953488372
$highValue = 33 / 43    # 0.76744186046511631...(rounds to .76744186046511631 at double precision)

# Switch-over column (1-based, A=1) per row: everything before this column
# gets $lowValue, everything from this column onward gets $highValue.
$switchCol = @{ 1 = 52; 2 = 51; 3 = 52; 4 = 50; 5 = 51 }   # AZ, AY, AZ, AX, AY

$acr = $wb.Worksheets.Item("ACR")
$numCols = 101  # columns A..CW
$numRows = 5

$data = New-Object 'object[,]' $numRows, $numCols
for ($r = 1; $r -le $numRows; $r++) {
    $sw = $switchCol[$r]
    for ($c = 1; $c -le $numCols; $c++) {
        if ($c -lt $sw) {
            $data[$r - 1, $c - 1] = $lowValue
        } else {
            $data[$r - 1, $c - 1] = $highValue
        }
    }
}
$acr.Range("A1:CW5").Value = $data

# FAR sheet: shift the FAR=1 / FAR=0 boundary at AX/AY in rows 1 and 3,
# and at AW/AX in row 4.
$far = $wb.Worksheets.Item("FAR")
$far.Range("AY1").Value = 1
$far.Range("AY3").Value = 1
$far.Range("AX4").Value = 0

# FRR sheet: the complementary boundary shift.
$frr = $wb.Worksheets.Item("FRR")
$frr.Range("AY1").Value = 0
$frr.Range("AY3").Value = 0
$frr.Range("AX4").Value = 1
